$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Nazwa testu) rows 6,7,8 in top-to-bottom order -> new shared strings 18,19,20
$ws.Range("C6").Value = "test_put_product_happy"
$ws.Range("C7").Value = "test_put_product_negative_valid_input_not_found"
$ws.Range("C8").Value = "test_put_product_negative_invalid_input_schema"

# Column B (Metoda) rows 6,7,8 -> new shared string 21 (reused for rows 7 and 8)
$ws.Range("B6").Value = "PUT"
$ws.Range("B7").Value = "PUT"
$ws.Range("B8").Value = "PUT"

# Column E (Krotki opis) bottom-to-top: row 8, row 7, row 6 -> new shared strings 22,23,24
$ws.Range("E8").Value = "sprawdza czy api poprawnie zwraca kod błędu przy próbie błędnej aktualizacji danych"
$ws.Range("E7").Value = "sprawdza czy api zwraca poprawny kod błędu przy próbie aktualizacji nieistniejącego produktu"
$ws.Range("E6").Value = "sprawdza czy api poprawnie zwraca kod oraz aktualizuje zasób"

# Column D (Kategoria) - reuses existing shared strings
$ws.Range("D6").Value = "happy tests"
$ws.Range("D7").Value = "negative testing with valid input"
$ws.Range("D8").Value = "negative testing with invalid input"

# Rows grow taller to fit the new wrapped text (matches the authored row heights)
$ws.Rows("6").RowHeight = 34
$ws.Rows("7").RowHeight = 34
$ws.Rows("8").RowHeight = 34

# Cursor ends up on E7 after filling the table in
$ws.Range("E7").Select()
